$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $ws.Range("E$row").Value = 1934
    $currentAge = $ws.Range("G$row").Value2
    $ws.Range("G$row").Value = $currentAge + 1
}
